# Refresh the crypto price/volume table (GitHub Actions bot update).
# Price (column D) and Volume(1h) (column E) are stored as plain text in the
# sheet, so numeric-looking price strings are written with a leading
# apostrophe to keep Excel from auto-converting them to numbers; the cell
# Style is then reset to "Normal" so that text-forcing quote-prefix marker
# doesn't leave a stray number-format style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.028.18'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.886.13'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''331.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '''0.4602'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.00%  '
$ws.Range('D8').Value = '''0.4066'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('D10').Value = '''0.07983'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('D11').Value = '''0.9919'
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Value = '''21.68'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').Value = '1.876.22'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').Value = '''5.912'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').Value = '''7.070'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.27%  '
$ws.Range('D16').Value = '''1.002'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '''88.51'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.39%  '
$ws.Range('D18').Value = '''0.00001029'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('D19').Value = '''0.06553'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').Value = '''17.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Value = '29.075.89'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = '''5.422'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('D24').Value = '''11.44'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').Value = '''2.214'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('D26').Value = '2.129.06'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '''157.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').Value = '''19.58'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.35%  '
$ws.Range('D29').Value = '''2.098'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.88%  '
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').Value = '''117.55'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.93%  '
$ws.Range('D32').Value = '''1.010'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').Value = '''0.09332'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '''1.409'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('D36').Value = '''5.279'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.82%  '
$ws.Range('D37').Value = '''0.06063'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').Value = '''0.02222'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').Value = '''8.303'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.43%  '
$ws.Range('D40').Value = '''1.176'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('D41').Value = '''1.001'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '''0.5786'
$ws.Range('D42').Style = "Normal"
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''10.14'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.20%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '''0.1822'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.16%  '
$ws.Range('D45').Value = '''1.260'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('D46').Value = '''0.07439'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('D47').Value = '''2.271'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.96%  '
$ws.Range('D48').Value = '''12.01'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('D49').Value = '''0.5444'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.45%  '
$ws.Range('E50').Value = '  -4.06%  '
$ws.Range('D51').Value = '''45.61'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +11.91%  '
